$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the General Dynamics (GD) contract term end date from
#    "Present" to "April 2019" so the term dates are consistent with the
#    actual end date ("Fixed term date for GD to be consistent").
#    The text reads: " 2018 <en dash> Present" -> " 2018 <en dash> April 2019"
# ---------------------------------------------------------------------------
$dash      = [char]8211
$oldSuffix = $dash + " Present"
$oldDate   = "2018 " + $oldSuffix
$newSuffix = $dash + " April 2019"
$newDate   = "2018 " + $newSuffix

$found = $d.Content.Find.Execute($oldDate, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $newDate, 2)
if (-not $found) {
    throw "Could not find the GD term-date text ('2018 - Present') to replace"
}

# ---------------------------------------------------------------------------
# 2. In the Solaris system-administration paragraph, the run boundary that
#    splits the word "services" moves over by one character - without any
#    change to the visible text ("...many services like NIS...").
#       before: run A ends "...many service" | run B starts "s like NIS..."
#       after : run A ends "...many serv"    | run B starts "ices like NIS..."
#    Both runs keep their original (identical) run formatting, so we rebuild
#    just the two runs in place via InsertXML rather than a plain text
#    replace (which the editor would otherwise coalesce into a single run).
# ---------------------------------------------------------------------------
$paraRange = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*many service*") {
        $paraRange = $p.Range
        break
    }
}
if ($paraRange -eq $null) {
    throw "Could not find the 'many services like NIS...' paragraph"
}

$fullText = $paraRange.Text
$marker = "many service"
$markerPos = $fullText.IndexOf($marker)
if ($markerPos -lt 0) {
    throw "Could not locate the run-split marker text"
}

# Text up to and including "...many serv" (run A, 1 char shorter than before)
$runAText = $fullText.Substring(0, $markerPos + $marker.Length - 1)
# Remainder of the paragraph text, starting with "ices like NIS..." (run B,
# 1 char longer than before). Trim the trailing paragraph-mark character
# that Range.Text includes at the end.
$remainder = $fullText.Substring($runAText.Length)
$remainder = $remainder.TrimEnd([char]13, [char]7)
$runBText = $remainder

$rPr = '<w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana" w:cs="Verdana"/><w:color w:val="auto"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>'

$newRunsXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">' + $runAText + '</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">' + $runBText + '</w:t></w:r>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Replace the runs (but not the paragraph mark) so the paragraph's own
# attributes/pPr are left completely untouched.
$runsRange = $d.Range($paraRange.Start, $paraRange.End - 1)
$runsRange.InsertXML($newRunsXml)
